$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new BOM rows (FD1,FD2,FD3 / Pipe), pushing the rest of the
#    table (and everything below it) down by two rows. We insert one row at
#    a time, right after the "D1" row (row 9), re-using row 9's formatting
#    (borders/fill/font/wrap) via a formats-only paste so the new rows look
#    just like the rest of the table body.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Insert()
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(12).Insert()
$ws.Range("A9:H9").Copy()
$ws.Range("A12:H12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Header block: BOM now generated from the PcbDoc (not just the PrjPcb),
#    so the title / "Source Data From" / report time change.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Bill of Materials for PCB Document [PCB_dHandlebar_Lateral_Left.PcbDoc]"
$ws.Range("D3").Value = "PCB_dHandlebar_Lateral_Left.PcbDoc"
$ws.Range("D4").Value = "PCB_dHandlebar_Lateral.PrjPcb"
$ws.Range("D6").Value = "11:34"

# ---------------------------------------------------------------------------
# 3. BOM rows 9-14 (Designator / LibRef / Description / Supplier 1 /
#    Supplier Part Number 1 / Quantity).
# ---------------------------------------------------------------------------

# Row 9 - D1 (unchanged designator/description, new supplier info)
$ws.Range("C9").Value = "D1"
$ws.Range("D9").Value = "LED RGB QLS6B-FKW-CNSNSF043"
$ws.Range("E9").Value = "LED, Superior, Transparente, RGB, Rojo, Verde, Azul, SMD, 120°, Redondo, R 30mA, V 30mA, B 30mA"
$ws.Range("F9").Value = "Mouser"
$ws.Range("G9").Value = "941-QLS6BFKWNSNSF043"
$ws.Range("H9").Value = 1

# Row 10 - new Fiducials row
$ws.Range("B10").Formula = "=ROW(B10) - ROW(`$B`$8)"
$ws.Range("C10").Value = "FD1, FD2, FD3"
$ws.Range("D10").Value = "FIDUCIAL"
$ws.Range("E10").Value = "Fiducial"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = 3

# Row 11 - Connector pads (was row 10)
$ws.Range("C11").Value = "P1, P2, P3, P4, P5, P6, P7, P8, P9"
$ws.Range("D11").Value = "CONNECTOR PAD"
$ws.Range("E11").Value = "Conector Pad"
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = 9

# Row 12 - new Pipe row
$ws.Range("B12").Formula = "=ROW(B12) - ROW(`$B`$8)"
$ws.Range("C12").Value = "Pipe"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = 1

# Row 13 - Resistors (was row 11), new supplier info
$ws.Range("C13").Value = "R1, R2, R3"
$ws.Range("D13").Value = "R 0805"
$ws.Range("E13").Value = "Resistencia SMD de Tipo Chip, Película Gruesa, 180 ohm, ± 1%, 500 mW, 0805 [Métrica 2012]"
$ws.Range("F13").Value = "Mouser"
$ws.Range("G13").Value = "667-ERJ-P06F1800V"
$ws.Range("H13").Value = 3

# Row 14 - Switch U1 (was row 12), new supplier info
$ws.Range("C14").Value = "U1"
$ws.Range("D14").Value = "KSC541J"
$ws.Range("E14").Value = "Interruptor Táctil, Sealed, KSC Series, Accionamiento Superior, Montaje Superficial, Botón Redondo"
$ws.Range("F14").Value = "Mouser"
$ws.Range("G14").Value = "611-KSC541JROHS"
$ws.Range("H14").Value = 1

# ---------------------------------------------------------------------------
# 4. Footer: total component count (was 14, now 18).
# ---------------------------------------------------------------------------
$ws.Range("H15").Value = "18"

# ---------------------------------------------------------------------------
# 5. Selection, matching what was captured when the workbook was re-saved.
# ---------------------------------------------------------------------------
$ws.Range("B9:H12").Select()
